$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grouped matches")

$ws.Cells.Item(2, 2).Value = "https://lod.dila.edu.tw/resource.php?id=A000089"
$ws.Cells.Item(2, 3).Value = "{'eft:siladharma'}"
$ws.Cells.Item(3, 2).Value = "P4CZ15137"
$ws.Cells.Item(3, 3).Value = "{'eft:kumarakalasa'}"
$ws.Cells.Item(4, 2).Value = "P8269"
$ws.Cells.Item(4, 3).Value = "{'eft:dgon-gling-rma'}"
$ws.Cells.Item(5, 2).Value = "P4CZ16819"
$ws.Cells.Item(5, 3).Value = "{'eft:sakyaprabha'}"
$ws.Cells.Item(6, 2).Value = "P0TMP092"
$ws.Cells.Item(6, 3).Value = "{'eft:anandasri-s-'}"
$ws.Cells.Item(7, 2).Value = "P4258"
$ws.Cells.Item(7, 3).Value = "{'eft:dpal-byor'}"
$ws.Cells.Item(8, 2).Value = "P8222"
$ws.Cells.Item(8, 3).Value = "{'eft:jnanasidhi', 'eft:jnanasiddhi'}"
$ws.Cells.Item(9, 2).Value = "P8266"
$ws.Cells.Item(9, 3).Value = "{'eft:ch-nyi-tsultrim', 'eft:dharmatasila'}"
$ws.Cells.Item(10, 2).Value = "P4259"
$ws.Cells.Item(10, 3).Value = "{'eft:ban-de-dpal-gyi-lhun-po', 'eft:palgyi-lh-npo', 'eft:dpal-gyi-lhun-po'}"
$ws.Cells.Item(11, 2).Value = "P00KG07267"
$ws.Cells.Item(11, 3).Value = "{'eft:sarvajnadeva', 'eft:sarvanyadeva'}"
$ws.Cells.Item(12, 2).Value = "?"
$ws.Cells.Item(12, 3).Value = "{'eft:sakyasena'}"
$ws.Cells.Item(13, 2).Value = "P3285"
$ws.Cells.Item(13, 3).Value = "{'eft:sakya-yesh-'}"
$ws.Cells.Item(14, 2).Value = "P2548"
$ws.Cells.Item(14, 3).Value = "{'eft:prajnavarman', 'eft:prajnavarma'}"
$ws.Cells.Item(15, 2).Value = "P8205"
$ws.Cells.Item(15, 3).Value = "{'eft:band-yesh-d-', 'eft:zhang-yesh-d-', 'eft:ye-shes-sde', 'eft:band-yesh-de', 'eft:yesh-d-ye-shes-sde-', 'eft:yesh-d-'}"
$ws.Cells.Item(16, 2).Value = "P8182"
$ws.Cells.Item(16, 3).Value = "{'eft:ban-de-dpal-brtsegs', 'eft:kawa-paltsek-under-the-name-paltsek-raksita-', 'eft:ska-ba-dpal-brtsegs', 'eft:dpal-brtsegs', 'eft:paltsek'}"
$ws.Cells.Item(17, 2).Value = "P4255"
$ws.Cells.Item(17, 3).Value = "{'eft:ye-shes-snying-po', 'eft:t-jnanagarbha', 'eft:yesh-nyingpo'}"
$ws.Cells.Item(18, 2).Value = "P3379"
$ws.Cells.Item(18, 3).Value = "{'eft:dipamkarasrijnana', 'eft:dipamkara-srijnana'}"
$ws.Cells.Item(19, 2).Value = "P8267"
$ws.Cells.Item(19, 3).Value = "{'eft:vijayasila'}"
$ws.Cells.Item(20, 2).Value = "P8260"
$ws.Cells.Item(20, 3).Value = "{'eft:dpal-dbyangs'}"
$ws.Cells.Item(21, 2).Value = "P1KG8854"
$ws.Cells.Item(21, 3).Value = "{'eft:srilendrabodhi', 'eft:surendrabodhi', 'eft:silendrabodhi'}"
$ws.Cells.Item(22, 2).Value = "P8171"
$ws.Cells.Item(22, 3).Value = "{'eft:dharmasribhadra'}"
$ws.Cells.Item(23, 2).Value = "P2637"
$ws.Cells.Item(23, 3).Value = "{'eft:trakpa-gyaltsen'}"
$ws.Cells.Item(24, 2).Value = "P8273"
$ws.Cells.Item(24, 3).Value = "{'eft:rinchen-tso', 'eft:rin-chen-tsho'}"
$ws.Cells.Item(25, 2).Value = "P8206"
$ws.Cells.Item(25, 3).Value = "{'eft:celu'}"
$ws.Cells.Item(26, 2).Value = "P8249"
$ws.Cells.Item(26, 3).Value = "{'eft:dharmakara'}"
$ws.Cells.Item(27, 2).Value = "P8210"
$ws.Cells.Item(27, 3).Value = "{'eft:danasila'}"
$ws.Cells.Item(28, 2).Value = "P8245"
$ws.Cells.Item(28, 3).Value = "{'eft:buddhakaravarma'}"
$ws.Cells.Item(29, 2).Value = "P4263"
$ws.Cells.Item(29, 3).Value = "{'eft:dge-ba-dpal'}"
$ws.Cells.Item(30, 2).Value = "P8209"
$ws.Cells.Item(30, 3).Value = "{'eft:jinamitra', 'eft:jinamitra-k-', 'eft:dzi-na-mi-tra-k-'}"
$ws.Cells.Item(31, 2).Value = "P2956"
$ws.Cells.Item(31, 3).Value = "{'eft:krsnapandita'}"
$ws.Cells.Item(32, 2).Value = "P8093"
$ws.Cells.Item(32, 3).Value = "{'eft:kamalagupta'}"
$ws.Cells.Item(33, 2).Value = "P0RK8"
$ws.Cells.Item(33, 3).Value = "{'eft:dharmapala'}"
$ws.Cells.Item(34, 2).Value = "P8228"
$ws.Cells.Item(34, 3).Value = "{'eft:surendrabodhi'}"
$ws.Cells.Item(35, 2).Value = "P8151"
$ws.Cells.Item(35, 3).Value = "{'eft:gayadhara'}"
$ws.Cells.Item(36, 2).Value = "P8261"
$ws.Cells.Item(36, 3).Value = "{'eft:munivarma', 'eft:munivarman'}"
$ws.Cells.Item(37, 2).Value = "P8211"
$ws.Cells.Item(37, 3).Value = "{'eft:vidyakaraprabha'}"
$ws.Cells.Item(38, 2).Value = "P0TMP104"
$ws.Cells.Item(38, 3).Value = "{'eft:punyasambhava'}"
$ws.Cells.Item(39, 2).Value = "P3709"
$ws.Cells.Item(39, 3).Value = "{'eft:phakpa-sherab'}"
$ws.Cells.Item(40, 2).Value = "P8219"
$ws.Cells.Item(40, 3).Value = "{'eft:visuddhasimha'}"
$ws.Cells.Item(41, 2).Value = "P8263"
$ws.Cells.Item(41, 3).Value = "{'eft:leki-d-'}"
$ws.Cells.Item(42, 2).Value = "P0TMP080"
$ws.Cells.Item(42, 3).Value = "{'eft:hwa-shang-zab-mo'}"
$ws.Cells.Item(43, 2).Value = "P3214"
$ws.Cells.Item(43, 3).Value = "{'eft:danasila'}"
$ws.Cells.Item(44, 2).Value = "P5651"
$ws.Cells.Item(44, 3).Value = "{'eft:pa-tshab-nyi-ma-grags'}"
$ws.Cells.Item(45, 2).Value = "P8213"
$ws.Cells.Item(45, 3).Value = "{'eft:vidyakarasimha', 'eft:t-vidyakarasimha'}"
$ws.Cells.Item(46, 2).Value = "P8220"
$ws.Cells.Item(46, 3).Value = "{'eft:devacandra'}"
$ws.Cells.Item(47, 2).Value = "P8265"
$ws.Cells.Item(47, 3).Value = "{'eft:ratnaraksita'}"
$ws.Cells.Item(48, 2).Value = "P4CZ16780"
$ws.Cells.Item(48, 3).Value = "{'eft:manjusrigarbha'}"
$ws.Cells.Item(49, 2).Value = "P4242"
$ws.Cells.Item(49, 3).Value = "{'eft:sherab-lekpa'}"
$ws.Cells.Item(50, 2).Value = "P8268"
$ws.Cells.Item(50, 3).Value = "{'eft:buddhaprabha'}"
$ws.Cells.Item(51, 2).Value = "P8217"
$ws.Cells.Item(51, 3).Value = "{'eft:jnanagarbha', 'eft:t-jnanagarbha'}"
$ws.Cells.Item(52, 2).Value = "P0TMPT007"
$ws.Cells.Item(52, 3).Value = "{'eft:rnam-par-mi-rtog-pa'}"
$ws.Cells.Item(53, 2).Value = "P0TMP098"
$ws.Cells.Item(53, 3).Value = "{'eft:jinavara'}"
$ws.Cells.Item(54, 2).Value = "P753"
$ws.Cells.Item(54, 3).Value = "{'eft:rin-chen-bzang-po'}"
$ws.Cells.Item(55, 2).Value = "P8183"
$ws.Cells.Item(55, 3).Value = "{'eft:klu-i-rgyal-mtshan', 'eft:cog-ro-klu-i-rgyal-mtshan'}"
$ws.Cells.Item(56, 2).Value = "P3456"
$ws.Cells.Item(56, 3).Value = "{'eft:tshul-khrims-rgyal-ba'}"
